$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: prime formatting for brand-new cells (format only, via copy/paste-special) ---
$ws.Range("B16:C16").Copy() | Out-Null
$ws.Range("B17:C17").PasteSpecial(-4122) | Out-Null
$ws.Range("A21:C21").Copy() | Out-Null
$ws.Range("A22:C22").PasteSpecial(-4122) | Out-Null
$ws.Range("A21:C21").Copy() | Out-Null
$ws.Range("A23:C23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Step 2: clear cells that become empty ---
$ws.Range("A13").ClearContents() | Out-Null
$ws.Range("A14").ClearContents() | Out-Null
$ws.Range("B19").ClearContents() | Out-Null
$ws.Range("C19").ClearContents() | Out-Null

# --- Step 3: set changed cell values (only cells whose content actually changes) ---
$ws.Range("B10").Value = 'Dar o embasamento dos conceitos elementares em química aos alunos, capacitando-os para o prosseguimento dos estudos nas disciplinas correlatas posteriores, principalmente quanto aos conceitos da estrutura atômica; das ligações químicas, forças intermoleculares e natureza dos compostos; da geometria das moléculas; das reações químicas em solução aquosa, tanto de dupla-troca como de oxi-redução; das propriedades do estado gasoso e das soluções e da estequiometria e cálculos em química, com ênfase em casos contendo reagentes limitantes, pureza de reagentes e rendimento de reação.'
$ws.Range("C10").Value = 'Dar o embasamento dos conceitos elementares em química aos alunos, capacitando-os para o prosseguimento dos estudos nas disciplinas correlatas posteriores, principalmente quanto aos conceitos da estrutura atômica; das ligações químicas, forças intermoleculares e natureza dos compostos; da geometria das moléculas; das reações químicas em solução aquosa, tanto de dupla-troca como de oxi-redução; das propriedades do estado gasoso e das soluções e da estequiometria e cálculos em química, com ênfase em casos contendo reagentes limitantes, pureza de reagentes e rendimento de reação.'
$ws.Range("B13").Value = '198273 - Domingos Savio Giordani'
$ws.Range("C13").Value = '198273 - Domingos Savio Giordani'
$ws.Range("B14").Value = '1506103 - Pedro Carlos de Oliveira'
$ws.Range("C14").Value = '1506103 - Pedro Carlos de Oliveira'
$ws.Range("A15").Value = 'Programa resumido:'
$ws.Range("B15").Value = 'Princípios elementares em química. Estrutura Atômica e Tabela Periódica. A Ligação Química. Natureza dos Compostos. Reações Químicas em Solução Aquosa. Gases. Soluções. Estequiometria e Cálculos em Química.'
$ws.Range("C15").Value = 'Princípios elementares em química. Estrutura Atômica e Tabela Periódica. A Ligação Química. Natureza dos Compostos. Reações Químicas em Solução Aquosa. Gases. Soluções. Estequiometria e Cálculos em Química.'
$ws.Range("A16").Value = 'Short syllabus:'
$ws.Range("B16").Value = 'Elementary principles of chemistry. Atomic structure and the Periodic Table. The Chemical Bonding. Nature of the compounds. Chemical Reactions in Aqueous Solution. Gases. Solutions. Stoichiometry calculations in chemistry.'
$ws.Range("C16").Value = 'Elementary principles of chemistry. Atomic structure and the Periodic Table. The Chemical Bonding. Nature of the compounds. Chemical Reactions in Aqueous Solution. Gases. Solutions. Stoichiometry calculations in chemistry.'
$ws.Range("A17").Value = 'Programa:'
$ws.Range("B17").Value = 'Princípios elementares em química: Sistemas de Unidades (Definição das Unidades mais usadas em Engenharia e transformações entre sistemas).Estrutura Atômica e Tabela Periódica: Natureza elétrica da matéria. A carga do elétron. O núcleo do átomo. Espectros de emissão e de absorção atômica. Configuração eletrônica dos elementos. Partículas Elementares. A Lei e a tabela Periódica.A Ligação Química: A ligação eletrovalente. A ligação covalente. Hibridação. Polaridade da ligação. Natureza dos Compostos: Ácidos e bases (Arrhenius, Bronsted-Lowry e Lewis). Forças intermoleculares.Reações Químicas em Solução Aquosa : Terminologia das soluções. Eletrólitos e não eletrólitos. Reações iônicas. Reações sem transferência de elétron e seu balanceamento. Preparação de sais inorgânicos (por dupla troca). Oxidação e redução. Número de oxidação. Reações de óxido redução. Métodos de balanceamento de reações de oxi-redução (Variação do Nox, via decomposição do agente oxidante, íon-elétron e pelo Potencial Padrão de Redução).Gases: Variáveis de estado. Lei combinada dos gases. Experiência de Torriceli. Teoria cinética dos gases. Gás ideal e real. Princípio de Avogadro.Soluções: Natureza das soluções. Dispersões coloidais e suspensões. Tipos de soluções. Unidades de concentração (Molaridade, fração molar, ppm, normalidade, molalidade). O processo de dissolução. Calor de dissolução. Solubilidade e temperatura.Estequiometria e Cálculos em Química : Cálculos baseados em equações químicas. Cálculos com reagentes limitantes e reagentes com pureza. Rendimento teórico e centesimal. Resolução de exercícios envolvendo estequiometria industrial.'
$ws.Range("C17").Value = 'Princípios elementares em química: Sistemas de Unidades (Definição das Unidades mais usadas em Engenharia e transformações entre sistemas).Estrutura Atômica e Tabela Periódica: Natureza elétrica da matéria. A carga do elétron. O núcleo do átomo. Espectros de emissão e de absorção atômica. Configuração eletrônica dos elementos. Partículas Elementares. A Lei e a tabela Periódica.A Ligação Química: A ligação eletrovalente. A ligação covalente. Hibridação. Polaridade da ligação. Natureza dos Compostos: Ácidos e bases (Arrhenius, Bronsted-Lowry e Lewis). Forças intermoleculares.Reações Químicas em Solução Aquosa : Terminologia das soluções. Eletrólitos e não eletrólitos. Reações iônicas. Reações sem transferência de elétron e seu balanceamento. Preparação de sais inorgânicos (por dupla troca). Oxidação e redução. Número de oxidação. Reações de óxido redução. Métodos de balanceamento de reações de oxi-redução (Variação do Nox, via decomposição do agente oxidante, íon-elétron e pelo Potencial Padrão de Redução).Gases: Variáveis de estado. Lei combinada dos gases. Experiência de Torriceli. Teoria cinética dos gases. Gás ideal e real. Princípio de Avogadro.Soluções: Natureza das soluções. Dispersões coloidais e suspensões. Tipos de soluções. Unidades de concentração (Molaridade, fração molar, ppm, normalidade, molalidade). O processo de dissolução. Calor de dissolução. Solubilidade e temperatura.Estequiometria e Cálculos em Química : Cálculos baseados em equações químicas. Cálculos com reagentes limitantes e reagentes com pureza. Rendimento teórico e centesimal. Resolução de exercícios envolvendo estequiometria industrial.'
$ws.Range("A18").Value = 'Syllabus:'
$ws.Range("B18").Value = 'Elementary principles of chemistry: Units Systems (Definition of the most used units in Engineering and transformations between systems).Atomic structure and the Periodic Table: electrical nature of matter. The electron charge. The nucleus of the atom. Emission spectra and atomic absorption. Electronic configuration of the elements. Elementary Particles. The Law and the Periodic Table.The Chemical Bonding: The ionic bonding. The covalent bond. Hybridization. Polarity of covalent bonding. Nature of the Compounds: Acids and bases (Arrhenius, Bronsted-Lowry and Lewis). Intermolecular forces.Chemical Reactions in Aqueous Solution: Terminology in Solutions. Electrolytes and non electrolytes. Ionic reactions. Reactions without electron transfer and its balancing. Preparation of inorganic salts (metathesis). Oxidation and reduction. Oxidation number. Redox reactions. Redox reactions balancing methods (Variation of Nox,  decomposition of the oxidizing agent, ion-electron and using the Standard Potential of Reduction).Gases: State variables. Combined gas law. Experience Torriceli. Kinetic theory of gases. Ideal and real gas. Avogadro''s Principle.Solutions: Nature of solutions. Colloidal dispersions and suspensions. Types of solutions. Concentration units (Molarity, mole fraction, ppm, normality, molality). The dissolution process. Heat dissolution. Solubility and temperature.Stoichiometric calculations in Chemistry: Calculations based on chemical equations. Calculations with limiting reagents and reagent purity. Theoretical and centesimal yields. Solving of exercises with industrial stoichiometric approach.'
$ws.Range("C18").Value = 'Elementary principles of chemistry: Units Systems (Definition of the most used units in Engineering and transformations between systems).Atomic structure and the Periodic Table: electrical nature of matter. The electron charge. The nucleus of the atom. Emission spectra and atomic absorption. Electronic configuration of the elements. Elementary Particles. The Law and the Periodic Table.The Chemical Bonding: The ionic bonding. The covalent bond. Hybridization. Polarity of covalent bonding. Nature of the Compounds: Acids and bases (Arrhenius, Bronsted-Lowry and Lewis). Intermolecular forces.Chemical Reactions in Aqueous Solution: Terminology in Solutions. Electrolytes and non electrolytes. Ionic reactions. Reactions without electron transfer and its balancing. Preparation of inorganic salts (metathesis). Oxidation and reduction. Oxidation number. Redox reactions. Redox reactions balancing methods (Variation of Nox,  decomposition of the oxidizing agent, ion-electron and using the Standard Potential of Reduction).Gases: State variables. Combined gas law. Experience Torriceli. Kinetic theory of gases. Ideal and real gas. Avogadro''s Principle.Solutions: Nature of solutions. Colloidal dispersions and suspensions. Types of solutions. Concentration units (Molarity, mole fraction, ppm, normality, molality). The dissolution process. Heat dissolution. Solubility and temperature.Stoichiometric calculations in Chemistry: Calculations based on chemical equations. Calculations with limiting reagents and reagent purity. Theoretical and centesimal yields. Solving of exercises with industrial stoichiometric approach.'
$ws.Range("A19").Value = 'Avaliação:'
$ws.Range("A20").Value = 'Método:'
$ws.Range("B20").Value = 'Duas provas escritas'
$ws.Range("C20").Value = 'Duas provas escritas'
$ws.Range("A21").Value = 'Critério:'
$ws.Range("B21").Value = 'A média para a primeira avaliação será calculada a partir das notas das duas provas, P1 e P2, segundo a fórmula: M1=(P1+2xP2)/3. Alunos com nota final igual ou superior a 5,0 estão aprovados; inferior a 5,0 e igual ou superior a 3,0 estão de recuperação; inferior a 3 estão reprovados.'
$ws.Range("C21").Value = 'A média para a primeira avaliação será calculada a partir das notas das duas provas, P1 e P2, segundo a fórmula: M1=(P1+2xP2)/3. Alunos com nota final igual ou superior a 5,0 estão aprovados; inferior a 5,0 e igual ou superior a 3,0 estão de recuperação; inferior a 3 estão reprovados.'
$ws.Range("A22").Value = 'Norma de recuperação:'
$ws.Range("B22").Value = 'A recuperação consistirá em uma prova envolvendo o assunto do semestre todo, à qual será atribuída nota NR. A média da segunda avaliação será calculada segundo a fórmula: M2=(M1+NR)/2. Alunos com nota M2 igual ou superior a 5,0 estarão aprovados, inferior a 5,0 estarão reprovados.'
$ws.Range("C22").Value = 'A recuperação consistirá em uma prova envolvendo o assunto do semestre todo, à qual será atribuída nota NR. A média da segunda avaliação será calculada segundo a fórmula: M2=(M1+NR)/2. Alunos com nota M2 igual ou superior a 5,0 estarão aprovados, inferior a 5,0 estarão reprovados.'
$ws.Range("A23").Value = 'Bibliografia:'
$ws.Range("B23").Value = 'BROWN, T.L. ET al. Química a ciência central. 9.ed. São Paulo: Pearson Prentice Hall, 2005-2007ATKINS, Peter., Princípios de Química, questionando a vida moderna e o meio ambiente. 3ª Ed. Porto Alegre: Editora Bookman, 2006BRADY, J ; HUMISTON, G.E. Química geral. Rio de Janeiro: Ed. Livros Técnicos Científicos, 1981CHANG, Raymond. Química geral: conceitos essenciais. 4.ed. AMGH Editora Ltda., 2010.RUSSEL, J.B. Química geral. São Paulo: MacGrall-Hill'
$ws.Range("C23").Value = 'BROWN, T.L. ET al. Química a ciência central. 9.ed. São Paulo: Pearson Prentice Hall, 2005-2007ATKINS, Peter., Princípios de Química, questionando a vida moderna e o meio ambiente. 3ª Ed. Porto Alegre: Editora Bookman, 2006BRADY, J ; HUMISTON, G.E. Química geral. Rio de Janeiro: Ed. Livros Técnicos Científicos, 1981CHANG, Raymond. Química geral: conceitos essenciais. 4.ed. AMGH Editora Ltda., 2010.RUSSEL, J.B. Química geral. São Paulo: MacGrall-Hill'

# --- Step 4: row heights ---
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 120
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 60
$ws.Rows.Item(23).RowHeight = 120
$ws.Rows.Item(13).AutoFit() | Out-Null
$ws.Rows.Item(14).AutoFit() | Out-Null
$ws.Rows.Item(19).AutoFit() | Out-Null
